$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# row 1 col 1: "53÷5=10, 3" -> "34÷6=5, 4"
$cell = $t.Cell(1, 1)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "34÷6=5, 4"

# row 1 col 2: "79÷5=15, 4" -> "42÷6=7, 0"
$cell = $t.Cell(1, 2)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "42÷6=7, 0"

# row 1 col 3: "43÷6=7, 1" -> "44÷9=4, 8"
$cell = $t.Cell(1, 3)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "44÷9=4, 8"

# row 1 col 4: "79÷5=15, 4" -> "44÷9=4, 8"
$cell = $t.Cell(1, 4)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "44÷9=4, 8"

# row 1 col 5: "84÷8=10, 4" -> "45÷9=5, 0"
$cell = $t.Cell(1, 5)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "45÷9=5, 0"

# row 5 col 1: "29÷7=4, 1" -> "71÷4=17, 3"
$cell = $t.Cell(5, 1)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "71÷4=17, 3"

# row 5 col 2: "13÷4=3, 1" -> "18÷6=3, 0"
$cell = $t.Cell(5, 2)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "18÷6=3, 0"

# row 5 col 3: "88÷7=12, 4" -> "62÷2=31, 0"
$cell = $t.Cell(5, 3)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "62÷2=31, 0"

# row 5 col 4: "55÷4=13, 3" -> "43÷8=5, 3"
$cell = $t.Cell(5, 4)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "43÷8=5, 3"

# row 5 col 5: "64÷7=9, 1" -> "50÷4=12, 2"
$cell = $t.Cell(5, 5)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "50÷4=12, 2"

# row 9 col 1: "21÷5=4, 1" -> "52÷7=7, 3"
$cell = $t.Cell(9, 1)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "52÷7=7, 3"

# row 9 col 2: "64÷3=21, 1" -> "42÷9=4, 6"
$cell = $t.Cell(9, 2)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "42÷9=4, 6"

# row 9 col 3: "45÷9=5, 0" -> "14÷9=1, 5"
$cell = $t.Cell(9, 3)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "14÷9=1, 5"

# row 9 col 4: "43÷9=4, 7" -> "67÷8=8, 3"
$cell = $t.Cell(9, 4)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "67÷8=8, 3"

# row 9 col 5: "76÷6=12, 4" -> "70÷6=11, 4"
$cell = $t.Cell(9, 5)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "70÷6=11, 4"

# row 13 col 1: "43÷9=4, 7" -> "88÷7=12, 4"
$cell = $t.Cell(13, 1)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "88÷7=12, 4"

# row 13 col 2: "30÷6=5, 0" -> "73÷8=9, 1"
$cell = $t.Cell(13, 2)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "73÷8=9, 1"

# row 13 col 3: "81÷9=9, 0" -> "64÷7=9, 1"
$cell = $t.Cell(13, 3)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "64÷7=9, 1"

# row 13 col 4: "77÷6=12, 5" -> "44÷8=5, 4"
$cell = $t.Cell(13, 4)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "44÷8=5, 4"

# row 13 col 5: "56÷7=8, 0" -> "31÷2=15, 1"
$cell = $t.Cell(13, 5)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "31÷2=15, 1"

# row 17 col 1: "85÷3=28, 1" -> "31÷3=10, 1"
$cell = $t.Cell(17, 1)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "31÷3=10, 1"

# row 17 col 2: "57÷4=14, 1" -> "49÷2=24, 1"
$cell = $t.Cell(17, 2)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "49÷2=24, 1"

# row 17 col 3: "69÷7=9, 6" -> "47÷3=15, 2"
$cell = $t.Cell(17, 3)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "47÷3=15, 2"

# row 17 col 4: "58÷4=14, 2" -> "45÷6=7, 3"
$cell = $t.Cell(17, 4)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "45÷6=7, 3"

# row 17 col 5: "95÷4=23, 3" -> "76÷2=38, 0"
$cell = $t.Cell(17, 5)
$cr = $cell.Range
$tr = $d.Range($cr.Start, $cr.End - 1)
$tr.Text = "76÷2=38, 0"

